# Applies the "test 2" Game of Thrones pool update:
#  - rewrites the Bonus Questions answer rows (26-31)
#  - adds three new bonus-question rows (32-34) with their merges
#  - updates the dimension / selection implicitly via the writes below

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 26: "Who dies first?" ----
$ws.Range("B26").Value = "Who dies first? (2pts) "
$ws.Range("G26").Value = "Tormund"

# ---- Row 27: "Who dies last?" ----
$ws.Range("B27").Value = "Who dies last? (2 pts)"
$ws.Range("G27").Value = "The Night King"

# ---- Row 28: "Winner of CLEGANEBOWL?" (text unchanged, answer changes) ----
$ws.Range("B28").Value = "Winner of CLEGANEBOWL? (2)"
$ws.Range("G28").Value = "Sandor Clegane"

# ---- Row 29: "Who kills the Night King?" (expanded text + new answer) ----
$ws.Range("B29").Value = "Who kills the Night King?…if killed (3) **If you think he lives put “safe” **"
$ws.Range("G29").Value = "John Snow"

# ---- Row 30: "Who kills Cercei?" (expanded text + new answer) ----
$ws.Range("B30").Value = "Who kills Cercei? … if killed (3) **If you think this bitch lives put “safe” **"
$ws.Range("G30").Value = "Jaime Lannister"

# ---- Row 31: now "Will we see an undead Ned Stark?" ----
$ws.Range("B31").Value = "Will we see an undead Ned Stark? (2) … yes or no answer"
$ws.Range("G31").Value = "No"

# ---- New row 32: Stark children over/under ----
$ws.Range("B32").Value = "Over / under how many Stark children survive: +1.5 (2)"
$ws.Range("G32").Value = "Under"
$ws.Range("B32:F32").Merge()

# ---- New row 33: Arya kill count over/under ----
$ws.Range("B33").Value = "Over / under the number of people Arya kills: +3.5 (2)"
$ws.Range("G33").Value = "Over"
$ws.Range("B33:F33").Merge()

# ---- New row 34: moved-down "Who wins the Game of Thrones???" question ----
$ws.Range("B34").Value = "Who wins the Game of Thrones??? (5) **If only one person guesses correctly they win 20% of pool and a prize**"
$ws.Range("G34").Value = "Daenerys Targaryen"
$ws.Range("B34:F34").Merge()

# Row heights for the new rows match the existing bonus-question rows
$ws.Rows.Item(32).RowHeight = 20
$ws.Rows.Item(33).RowHeight = 20
$ws.Rows.Item(34).RowHeight = 20

# Selection ends on the newly-relocated "who wins" question, matching the
# author's final cursor position after restructuring the bonus section.
$ws.Range("B31:F31").Select()
